$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column at M (old column M "check/pass" data shifts right to N)
$ws.Columns.Item(13).Insert() | Out-Null

# Populate the new "target drive" / "/t" option column header cells
$ws.Range("M1").Value = "target drive"
$ws.Range("M2").Value = "/t"

# Size the new column to fit its short header content
$ws.Columns.Item(13).ColumnWidth = 4.86

# "set start up" action no longer applies the "default" option
$ws.Range("G9").Clear() | Out-Null

# Update the active selection shown when the file is reopened
$ws.Range("C5:C8").Select() | Out-Null
